# Update the "COW reports" path entry on the include sheet to point at the
# new "#TEST IMAGES" folder (per commit: "Change imgeprocess folder list").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("include")
$ws.Range("B2").Value = "C:\Temp\#TEST IMAGES\"
